# Natmi following Dr Hou advice
# Rewrite the LR-pairs data to a full 3x3 Sending-cluster x Target-cluster
# matrix (ECs, FAPs, sCs) for the Slit2-Robo1 pair, extending rows 2-5 up
# to rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @{}
$data[2]  = @("ECs", "Slit2", "Robo1", "ECs", 2, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662, 2, 0.6666666666666666, 0.1757713333333334, 0.5273140000000001, 0.009444264308298454, 0.009444264308298454, 0.005876621577777779, 0.05288959420000001, 0.0001655824676354638, 0.0001655824676354637)
$data[3]  = @("ECs", "Slit2", "Robo1", "FAPs", 2, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662, 3, 1, 16.98312366666667, 50.949371, 0.9125100529581165, 0.9125100529581165, 0.567802434588889, 5.1102219113, 0.01599866981467349, 0.01599866981467349)
$data[4]  = @("ECs", "Slit2", "Robo1", "sCs", 2, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662, 3, 1, 1.452542333333333, 4.357627, 0.07804568273358503, 0.07804568273358505, 0.04856333201111111, 0.4370699881, 0.001368343400127671, 0.001368343400127671)
$data[5]  = @("FAPs", "Slit2", "Robo1", "ECs", 3, 1, 0.37963, 1.13889, 0.1990797397484571, 0.1990797397484571, 2, 0.6666666666666666, 0.1757713333333334, 0.5273140000000001, 0.009444264308298454, 0.009444264308298454, 0.06672807127333334, 0.6005526414600001, 0.001880161680611698, 0.001880161680611698)
$data[6]  = @("FAPs", "Slit2", "Robo1", "FAPs", 3, 1, 0.37963, 1.13889, 0.1990797397484571, 0.1990797397484571, 3, 1, 16.98312366666667, 50.949371, 0.9125100529581165, 0.9125100529581165, 6.447303237576667, 58.02572913818999, 0.1816622638607526, 0.1816622638607526)
$data[7]  = @("FAPs", "Slit2", "Robo1", "sCs", 3, 1, 0.37963, 1.13889, 0.1990797397484571, 0.1990797397484571, 3, 1, 1.452542333333333, 4.357627, 0.07804568273358503, 0.07804568273358505, 0.5514286460033333, 4.962857814029999, 0.01553731420709276, 0.01553731420709276)
$data[8]  = @("sCs", "Slit2", "Robo1", "ECs", 3, 1, 1.493861, 4.481583000000001, 0.7833876645691064, 0.7833876645691064, 2, 0.6666666666666666, 0.1757713333333334, 0.5273140000000001, 0.009444264308298454, 0.009444264308298454, 0.2625779397846667, 2.363201458062, 0.007398520160051292, 0.007398520160051292)
$data[9]  = @("sCs", "Slit2", "Robo1", "FAPs", 3, 1, 1.493861, 4.481583000000001, 0.7833876645691064, 0.7833876645691064, 3, 1, 16.98312366666667, 50.949371, 0.9125100529581165, 0.9125100529581165, 25.37042610381034, 228.333834934293, 0.7148491192826905, 0.7148491192826905)
$data[10] = @("sCs", "Slit2", "Robo1", "sCs", 3, 1, 1.493861, 4.481583000000001, 0.7833876645691064, 0.7833876645691064, 3, 1, 1.452542333333333, 4.357627, 0.07804568273358503, 0.07804568273358505, 2.169896342615667, 19.529067083541, 0.06114002512636461, 0.06114002512636461)

foreach ($r in 2..10) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $rowVals[$i]
    }
}
